# PCB rev2-230718 Update to Use throughhole audio jacks (C2939180)
#
# The BOM lists four connectors (Cassette in J6, Printer J8, Audio out J3,
# Cassette out J5) that used the SMT audio jack part C2884998. This part is
# being swapped for the throughhole equivalent C2939180. The footprint
# (PJ-320B-SMT) itself is unchanged - only the JLCPCB part number column
# changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldPart = "C2884998"
$newPart = "C2939180"

# Column D holds "JLCPCB Part #(optional)"; update every row whose part
# number currently matches the old SMT jack part number.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -eq $oldPart) {
        $cell.Value = $newPart
    }
}

# Reflect that the user's active selection ended up on D22 after the edit.
$ws.Range("D22").Select()
